$d = $word.ActiveDocument

# Locate the paragraph that anchors the end of the region to replace
# (the last paragraph of the original content, containing the
# "Back to Home" link) so the new HTML-document skeleton is swapped
# in for paragraphs 1..N while the trailing "</body></html>" wrapper
# paragraphs (and the final section break) are left untouched.
$count = $d.Paragraphs.Count
$endIndex = -1
for ($i = 1; $i -le $count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*Back to Home*") {
        $endIndex = $i
        break
    }
}
if ($endIndex -eq -1) {
    throw "Could not locate the 'Back to Home' paragraph to anchor the replacement."
}

$endPara = $d.Paragraphs($endIndex)
$rng = $d.Range(0, $endPara.Range.End)

$newBody = @'
<w:p>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>&lt;!DOCTYPE</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t xml:space="preserve"> html&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t>&lt;</w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>html</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t>&lt;</w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>head</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">    &lt;</w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>title&gt;</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>Books - E-Library&lt;/title&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">    &lt;link </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>rel</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>="</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>stylesheet</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve">" </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>href</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>="style.css"&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t>&lt;/head&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t>&lt;</w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>body</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>&gt;</w:t></w:r>
</w:p>
<w:p/>
<w:p>
  <w:r><w:t>&lt;</w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>header</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">    &lt;h1&gt;</w:t></w:r>
  <w:r>
    <w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/></w:rPr>
    <w:t>📚</w:t>
  </w:r>
  <w:r><w:t xml:space="preserve"> Available Books&lt;/h1&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">    &lt;</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>nav</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">        &lt;a </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>href</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>="index.html"&gt;Home&lt;/a&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">        &lt;a </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>href</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>="books.html"&gt;Books&lt;/a&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">    &lt;/</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>nav</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t>&lt;/header&gt;</w:t></w:r>
</w:p>
<w:p/>
<w:p>
  <w:r><w:t>&lt;</w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>section</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">    &lt;</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>ul</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">        &lt;</w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>li&gt;</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>HTML Basics&lt;/li&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">        &lt;</w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>li&gt;</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>CSS Guide&lt;/li&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">        &lt;</w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>li&gt;</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>JavaScript Introduction&lt;/li&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">        &lt;</w:t></w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>li&gt;</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>Python Programming&lt;/li&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">    &lt;/</w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>ul</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t>&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t>&lt;/section&gt;</w:t></w:r>
</w:p>
<w:p/>
<w:p>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>&lt;</w:t>
  </w:r>
  <w:proofErr w:type="gramStart"/>
  <w:r><w:t>footer</w:t></w:r>
  <w:proofErr w:type="gramEnd"/>
  <w:r><w:t>&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t xml:space="preserve">    &lt;p&gt;© 2026 My E-Library&lt;/p&gt;</w:t></w:r>
</w:p>
<w:p>
  <w:r><w:t>&lt;/footer&gt;</w:t></w:r>
</w:p>
'@

$rng.InsertXML($newBody)

Write-Output "Paragraphs after edit: $($d.Paragraphs.Count)"
